$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff changes cell B11 on the "Rules" sheet from the text "R40" to the
# text "1" - a brand new shared-string entry is appended to the workbook's
# string table and B11 is repointed at it (cell formatting / style untouched).
#
# A plain `$ws.Range("B11").Value = "1"` would be auto-detected as a NUMBER
# (matching how Excel treats typed numeric-looking input), which is not what
# the diff shows (B11 keeps referencing a shared string). To force the value
# to be stored as text without touching B11's existing style, stage the text
# in a scratch cell far away from the used range, copy it, and paste only the
# value (not formats) into B11 - this is the standard Excel technique for
# moving a literal string into a cell while leaving its number format alone.
$scratch = $ws.Range("Z20")
$scratch.Formula = '="1"'
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Remove the scratch row completely again so the sheet's used range / row
# and column structure is left exactly as it was before this script ran.
$scratch.EntireRow.Delete()

$wb.Save()
